$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at B so that existing custom column widths
# (originally at B, D, F, H, I, J) shift right to C, E, G, I, J, K -
# this lines most of them up with their new target positions.
$ws.Columns("B:B").Insert()

# Set new header values. Order matters for shared-string table layout:
# introduce brand-new strings (title_ar, description_ar, sku_number) in the
# same order they first appear in the target file, before touching the
# other (pre-existing) strings.
$ws.Range("D1").Value = "title_ar"
$ws.Range("G1").Value = "description_ar"
$ws.Range("B1").Value = "sku_number"

$ws.Range("A1").Value = "category"
$ws.Range("C1").Value = "title"
$ws.Range("E1").Value = "sale_price"
$ws.Range("F1").Value = "description"
$ws.Range("H1").Value = "purchase_price"
$ws.Range("I1").Value = "discount"
$ws.Range("J1").Value = "add_stock"
$ws.Range("K1").Value = "published"
$ws.Range("L1").Value = "images"

# Clear the old trailing columns that are no longer part of the table.
$ws.Range("M1:P1").ClearContents()

# Column widths (runtime quantizes ColumnWidth to 1/6-character pixel
# steps internally, so values below are chosen to land as close as
# possible on the desired stored widths).
$ws.Columns.Item(2).ColumnWidth = 10.0                 # B -> stored ~10.8333 (target 10.88671875)
$ws.Columns.Item(9).ColumnWidth = 8.166666666666666    # I -> stored 9        (target 9)
$ws.Columns.Item(10).ColumnWidth = 12.333333333333334  # J -> stored ~13.1667 (target 13.21875)
$ws.Columns.Item(8).ColumnWidth = 14.833333333333334   # H -> stored ~15.6667 (target 15.6640625)

# Selection and page setup
$ws.Range("E8").Select()
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$wb.Save()
